$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.589741468429565
$ws.Range("B1").Value = 4.617850303649902
$ws.Range("C1").Value = 6.742213726043701
$ws.Range("D1").Value = 7.840402126312256
$ws.Range("E1").Value = 5.075840950012207
